# Weekly update: insert the latest week's "Apio" (celery) price data
# (date 2022-07-06 / serial 44748) as a new Primera/Segunda row pair at
# the top of the Terminal La Palmera de La Serena data block (row 377).
# Inserting there pushes every subsequent row down by two, which both
# re-numbers the existing rows AND carries the former last pair
# (old rows 431-432, date 2021-11-19) down to the new rows 433-434 --
# exactly matching the target diff. No other cells need to be touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 377 downward by two, creating two fresh blank rows at 377-378.
$ws.Rows.Item(377).Resize(2).Insert()

# --- Row 377: Apio, Americana (o), Primera ---
$ws.Range("A377").Value = 8
$ws.Range("B377").Value = "Terminal La Palmera de La Serena"
$ws.Range("C377").Value = "Coquimbo"
$ws.Range("D377").Value = 44748
$ws.Range("E377").Value = 4
$ws.Range("F377").Value = 100112017
$ws.Range("G377").Value = "Apio"
$ws.Range("H377").Value = "Americana (o)"
$ws.Range("I377").Value = "Primera"
$ws.Range("J377").Value = 2000
$ws.Range("K377").Value = 8000
$ws.Range("L377").Value = 9000
$ws.Range("M377").Value = 8500
$ws.Range("N377").Value = "`$/docena de matas"
$ws.Range("O377").Value = "Provincia del Elquí"
$ws.Range("P377").Value = 1417
$ws.Range("Q377").Value = 6
$ws.Range("R377").Value = "Hortaliza"

# --- Row 378: Apio, Americana (o), Segunda ---
$ws.Range("A378").Value = 8
$ws.Range("B378").Value = "Terminal La Palmera de La Serena"
$ws.Range("C378").Value = "Coquimbo"
$ws.Range("D378").Value = 44748
$ws.Range("E378").Value = 4
$ws.Range("F378").Value = 100112017
$ws.Range("G378").Value = "Apio"
$ws.Range("H378").Value = "Americana (o)"
$ws.Range("I378").Value = "Segunda"
$ws.Range("J378").Value = 1200
$ws.Range("K378").Value = 7000
$ws.Range("L378").Value = 7500
$ws.Range("M378").Value = 7250
$ws.Range("N378").Value = "`$/docena de matas"
$ws.Range("O378").Value = "Provincia del Elquí"
$ws.Range("P378").Value = 1208
$ws.Range("Q378").Value = 6
$ws.Range("R378").Value = "Hortaliza"
